# Apply "Taking changes to local" edits to the "Test Cases" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Update the Results column (E) for the affected test case rows.
$ws.Range("E2").Value  = "SKIP"
$ws.Range("E21").Value = "PASS"
$ws.Range("E84").Value = "SKIP"
$ws.Range("E85").Value = "SKIP"
$ws.Range("E86").Value = "SKIP"
$ws.Range("E87").Value = "SKIP"
$ws.Range("E88").Value = "SKIP"
$ws.Range("E89").Value = "SKIP"

# Update the sheet view: scroll back to the top and select D2:D89.
$ws.Range("D2:D89").Select()
